# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures (VALOR MORA, Cant. Trabajadores, Cant. Periodos) ---
$ws.Range("E11").Value = 413352
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 6

# --- Make room for 3 new worker rows right after the current last data row (20) ---
$ws.Rows("21:23").Insert()

# The "closing" bottom-border look currently still sitting on row 20 needs to move
# down to the new last row (23) of the table.
$ws.Range("B20:J20").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Rows 20, 21 and 22 become ordinary interior rows, matching row 19's look.
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J22").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill in the new worker rows ---
# Row 21: JOSE RAMON PARRA TALAIGUA
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "73187053"
$ws.Range("D21").Value = "JOSE RAMON PARRA TALAIGUA"
$ws.Range("E21").Value = "2509"
$ws.Range("F21").Value = 62632
$ws.Range("G21").Value = 1565800

# Row 22: JORGE EMIRO PADILLA MORALES
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047391662"
$ws.Range("D22").Value = "JORGE EMIRO PADILLA MORALES"
$ws.Range("E22").Value = "2509"
$ws.Range("F22").Value = 62632
$ws.Range("G22").Value = 1565800

# Row 23: DUVAN DE JESUS CASTILLA LORDUY (new period 2509), now the closing row
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1047418083"
$ws.Range("D23").Value = "DUVAN DE JESUS CASTILLA LORDUY"
$ws.Range("E23").Value = "2509"
$ws.Range("F23").Value = 57200
$ws.Range("G23").Value = 1430000

# Center the "Periodo Mora" column across every data row, new style added in this edit.
$ws.Range("E16:E23").HorizontalAlignment = -4108   # xlCenter

Write-Host "Edit applied"
